# GF$21841 add note about coordination with Argonaut for subscription-based use
# cases and remove subscriptions
#
# The Subscription resource/profile is being pulled back out of this
# Consumer Client CapabilityStatement. Rather than delete the rows outright,
# the two cells that identify the Subscription resource/profile are
# commented-out in place by prefixing their text with "!" (the convention
# this spreadsheet uses to mark a row as excluded from the generated
# CapabilityStatement while still preserving the row for reference).

$wb = $excel.ActiveWorkbook

# "resources" sheet: row 2, column A holds the resource type ("Subscription").
# Prefix it with "!" so it's excluded from the generated CapabilityStatement.
$wsResources = $wb.Worksheets.Item("resources")
$wsResources.Range("A2").Value = "!Subscription"

# "profiles" sheet: row 10, column A holds the Subscription profile's
# canonical URL. Prefix it with "!" too.
$wsProfiles = $wb.Worksheets.Item("profiles")
$wsProfiles.Range("A10").Value = "!http://hl7.org/fhir/us/davinci-deqm/STU3/StructureDefinition/subscription-deqm"

# Leave a trail of the user's navigation while making this edit, matching
# where they ended up looking (profiles -> resources -> igs, with igs as the
# final active tab/selection).
[void]$wsProfiles.Range("A14").Select()
[void]$wsResources.Range("A15").Select()

$wsIgs = $wb.Worksheets.Item("igs")
[void]$wsIgs.Range("B21").Select()
